$d = $word.ActiveDocument

function Insert-BodyXml($range, $bodyXml) {
    $pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    [void]$range.InsertXML($pkg)
}

# --- Paragraph 1: merge the leading tab run into the following text run ---
$para1Xml = '<w:p><w:r><w:tab/><w:t xml:space="preserve">My favorite genre of games to play is racing games, so I decided to bring out my Xbox 360 and play Need for Speed Most Wanted 2012. </w:t></w:r><w:r><w:t>Need for Speed is my favorite game series. It sparked my love for cars and shaped my childhood. I can proudly say I’ve owned and played every NFS title in the series.</w:t></w:r></w:p>'
Insert-BodyXml $d.Paragraphs(1).Range $para1Xml

# --- Paragraph 3 (the one holding the _GoBack bookmark): replace with full text ---
$para3Xml = '<w:p><w:r><w:tab/></w:r><w:r><w:t xml:space="preserve">The object of the game is to complete races </w:t></w:r><w:r><w:t>to</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> unlock cars and upgrade parts. As you earn Speed Points, you can challenge the 10 Blacklist racers. They have special cars that you can </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>win</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> and their races are harder. The player can also engage in police pursuits to earn Speed Points. The only thing I really don’t like about this game is that you don’t really get punished for losing. If you lose a race, you can just restart. If you are caught by the police, you just lose the points you earned during the chase. </w:t></w:r></w:p>'
Insert-BodyXml $d.Paragraphs(3).Range $para3Xml

# --- Insert the new paragraphs after paragraph 3, before the final paragraph mark ---
$insertPos = $d.Paragraphs(3).Range.End
$insertRange = $d.Range($insertPos, $insertPos)
$newParasXml = '<w:p/><w:p><w:pPr><w:ind w:firstLine="709"/></w:pPr><w:r><w:t xml:space="preserve">The interesting part about the game is that to unlock the non-Blacklist cars in the game, you </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>have to</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> drive around and find them. </w:t></w:r><w:r><w:t xml:space="preserve">Each car can be found in 3 locations, so it’s fun to search for every car. Once you get a car, you need to complete certain challenges to unlock performance upgrades. These range from tires, nitrous, chassis, and a couple of others. Once you equip a certain part, you can use that part to unlock a Pro version of the part. These upgrades are one of the core mechanics of the game. The parts the player chooses all affect how the cars drive. For example, if you are driving a Ford Raptor, you can equip and Impact Protection Body and a Heavy Chassis to essentially create a tank. You’ll be able to ram cops and traffic without hindering your speed, however the cost is that the truck is slow. On the other hand, you can choose a Koenigsegg </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Agera</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> R which is the fastest car in the game. If you equip it with a Lightweight Chassis, Aero Body, and Track Tires, you can reach speeds in excess of 260 mph. The cost is that the car is light, and therefore cops can easily push you around and make you wreck.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="709"/></w:pPr></w:p><w:p><w:r><w:tab/><w:t>The</w:t></w:r><w:r><w:t xml:space="preserve"> player’s choices affect more than police chases. Certain races also require certain parts and cars. For example, a Lamborghini is fast, but it’s not a good choice for an </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>offroad</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> race and you are guaranteed not to win. You need to choose an all-wheel-drive car with </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Offroad</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Tires in order to win. Similarly, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Offroad</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Tires won’t help you win a race on tarmac. They aren’t as grippy as Track tires, so you’ll need to equip those along with Short Gears to accelerate quickly.</w:t></w:r></w:p><w:p/><w:p><w:r><w:tab/><w:t>Need for Speed Most Wanted 2012 doesn’t have as many in-depth choices compared to a game like Forza. But the choices that it does let the player make impact the gameplay in a meaningful way. I still thoroughly enjoy the game to this day.</w:t></w:r></w:p>'
Insert-BodyXml $insertRange $newParasXml

Write-Output $d.Paragraphs.Count
Write-Output $d.Content.Text
